$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), copying the formatting
# already used by the other header cells (e.g. H1) so the new columns
# match the existing header style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I/J columns, rows 2-12
$data = @(
    @(16, 16),
    @(6, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(4, 5),
    @(7, 8),
    @(6, 7),
    @(7, 8),
    @(10, 10),
    @(8, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
